# Update the course-link textbox (shape 3, "標題 1" text box) on slides 2, 3
# and 4: replace the Udemy course URL text with the new GitHub repo URL and
# drop the cached normAutofit fontScale (PowerPoint recomputes it itself,
# i.e. the shrink-to-fit percentage is cleared back to "no cached scale").

$p = $ppt.ActivePresentation

$newText = "https://github.com/peterhchen/900_MEAN_Proj"
$slideIndexes = @(2, 3, 4)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(3)
    $shape.TextFrame.TextRange.Text = $newText

    # Re-assert "shrink text on overflow" autofit so the stored normAutofit
    # element no longer carries a stale fontScale attribute.
    $shape.TextFrame.AutoSize = 2
}
